$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: refresh timestamp value (floating-point serial precision) ---
$ws.Range('A17').Value = 45048.75874868056

# --- Append new chat-log rows 18-44 ---
$chatRows = @(
    [PSCustomObject]@{ Row = 18; Ts = 45055.03911534722; User = 'Hi'; Bot = 'Good afternoon! How can I assist you?' },
    [PSCustomObject]@{ Row = 19; Ts = 45055.03966903935; User = 'Hi'; Bot = 'Greetings! How may I assist you?' },
    [PSCustomObject]@{ Row = 20; Ts = 45055.03969355324; User = 'my name is matthew'; Bot = 'Hello matthew. How can I help you?' },
    [PSCustomObject]@{ Row = 21; Ts = 45055.04034922454; User = 'Hi'; Bot = 'I''m doing great, thanks for asking! How can I assist you?' },
    [PSCustomObject]@{ Row = 22; Ts = 45055.04038950231; User = 'my name is matthew'; Bot = 'Hello matthew. How can I help you?' },
    [PSCustomObject]@{ Row = 23; Ts = 45055.04058368056; User = 'hello'; Bot = 'I''m doing great, thanks for asking! How can I assist you?' },
    [PSCustomObject]@{ Row = 24; Ts = 45055.04172206019; User = 'Can you help me?'; Bot = 'It''s hard to say without more information about your specific needs. Could you please provide more details?' },
    [PSCustomObject]@{ Row = 25; Ts = 45055.04176005787; User = 'my name is matthew'; Bot = 'Hello Matthew. How can I help you?' },
    [PSCustomObject]@{ Row = 26; Ts = 45055.04341940973; User = 'Hi'; Bot = 'Good morning! How can I help you today?' },
    [PSCustomObject]@{ Row = 27; Ts = 45055.04344451389; User = 'my name is matthew'; Bot = 'Hello Matthew. How can I help you?' },
    [PSCustomObject]@{ Row = 28; Ts = 45055.04355265047; User = 'Im looking for ransomware'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 29; Ts = 45055.04363553241; User = 'looking for ransomeware'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 30; Ts = 45055.04375228009; User = 'looking for ransomware'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 31; Ts = 45055.04381690973; User = 'How can you help me?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 32; Ts = 45055.04406925926; User = 'Have a nice day'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 33; Ts = 45055.04434280092; User = 'Hi'; Bot = 'Good afternoon! How can I assist you?' },
    [PSCustomObject]@{ Row = 34; Ts = 45055.04440151621; User = 'How can you help me?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 35; Ts = 45055.04458592593; User = 'I need some assistance'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 36; Ts = 45055.04481875; User = 'Can you help me choose?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 37; Ts = 45055.04500344907; User = 'What ransomware do you have?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 38; Ts = 45055.04504543982; User = 'What ransomware do you have?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 39; Ts = 45055.04507420139; User = 'How can you help me?'; Bot = 'Can you please clarify your question?' },
    [PSCustomObject]@{ Row = 40; Ts = 45055.04510114584; User = 'Hi'; Bot = 'Howdy! How can I help you today?' },
    [PSCustomObject]@{ Row = 41; Ts = 45055.04565913195; User = 'Hi'; Bot = 'Good afternoon! How can I assist you?' },
    [PSCustomObject]@{ Row = 42; Ts = 45055.04568487268; User = 'Can you help me?'; Bot = 'We offer a range of ransomware solutions to meet the unique needs of your organization. Some of our most popular options include Lockdown, Hive, and Beacon. Would you like me to provide you with more information about these solutions?' },
    [PSCustomObject]@{ Row = 43; Ts = 45055.04575633102; User = 'looking for ransomware'; Bot = 'There are options for large or small companies. Which do you prefer?.' },
    [PSCustomObject]@{ Row = 44; Ts = 45055.04609741081; User = 'How can you help me?'; Bot = 'I would be happy to make a recommendation, but I need more information about what you''re looking for.' }
)

foreach ($r in $chatRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Ts
    $ws.Cells.Item($r.Row, 1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
    $ws.Cells.Item($r.Row, 2).Value = $r.User
    $ws.Cells.Item($r.Row, 3).Value = $r.Bot
}

